$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.318.28'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '1.800.21'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.12'
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.574'
$ws.Range('E6').Value = '  +3.69%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '35.84'
$ws.Range('E8').Value = '  +9.90%  '
$ws.Range('E9').Value = '  +1.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0691'
$ws.Range('E10').Value = '  +0.54%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0965'
$ws.Range('E11').Value = '  +2.18%  '
$ws.Range('D12').Value = '2.061.53'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.58'
$ws.Range('E13').Value = '  +5.12%  '
$ws.Range('D14').Value = '1.815.09'
$ws.Range('E14').Value = '  +1.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.641'
$ws.Range('E15').Value = '  +1.42%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.48'
$ws.Range('E16').Value = '  +4.80%  '
$ws.Range('D17').Value = '34.334.31'
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.90'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.42'
$ws.Range('E19').Value = '  +0.30%  '
$ws.Range('D20').Value = '0.0₃0793'
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.53'
$ws.Range('E21').Value = '  +2.47%  '
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.15'
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('E24').Value = '  +3.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.78'
$ws.Range('E25').Value = '  +3.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.93'
$ws.Range('E26').Value = '  +8.95%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.82'
$ws.Range('E27').Value = '  +2.13%  '
$ws.Range('E28').Value = '  +2.00%  '
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.99'
$ws.Range('E30').Value = '  +0.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0529'
$ws.Range('E31').Value = '  +1.06%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.24'
$ws.Range('E32').Value = '  +1.17%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.82'
$ws.Range('E33').Value = '  +0.77%  '
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('D35').Value = '1.397.42'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.668'
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('E37').Value = '  -4.72%  '
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0189'
$ws.Range('E39').Value = '  -0.21%  '
$ws.Range('E40').Value = '  +10.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.959'
$ws.Range('E41').Value = '  +2.87%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '82.42'
$ws.Range('E42').Value = '  -3.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.83'
$ws.Range('E43').Value = '  -0.28%  '
$ws.Range('E44').Value = '  +0.52%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.40'
$ws.Range('E45').Value = '  -2.36%  '
$ws.Range('E46').Value = '  -3.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.03'
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('D48').Value = '1.961.61'
$ws.Range('E48').Value = '  +0.56%  '
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  +0.07%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '103.99'
$ws.Range('E50').Value = '  -0.73%  '
$ws.Range('E51').Value = '  +0.25%  '
